$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC_UserRegistration")

$ws.Range("L2").Value = 'VerifyTitle: City Market Norwalk - Online Grocery Supermarket with Home Delivery'
$ws.Range("L3").Value = 'VerifyElement: null'
$ws.Range("L4").Value = 'Click: null'
$ws.Range("L5").Value = 'VerifyElement: null'
$ws.Range("L6").Value = 'VerifyElement: null'
$ws.Range("L7").Value = 'Click: null'
$ws.Range("L8").Value = 'VerifyText: Create New Account'
$ws.Range("L9").Value = 'VerifyElement: null'
$ws.Range("L10").Value = 'VerifyElement: null'
$ws.Range("L11").Value = 'VerifyElement: null'
$ws.Range("L12").Value = 'VerifyElement: null'
$ws.Range("L13").Value = 'VerifyElement: null'
$ws.Range("L14").Value = 'VerifyElement: null'
$ws.Range("L15").Value = 'VerifyElement: null'
$ws.Range("L16").Value = 'VerifyElement: null'
$ws.Range("L17").Value = 'VerifyElement: null'
$ws.Range("L18").Value = 'SetText: Akash'
$ws.Range("L19").Value = 'SetText: sangal'
$ws.Range("L20").Value = 'ClearText: null'
$ws.Range("L21").Value = 'ClearText: null'
$ws.Range("L22").Value = 'ClearText: null'
$ws.Range("L23").Value = 'Click: null'
$ws.Range("L24").Value = 'VerifyText: Please enter your email address.'
$ws.Range("L25").Value = 'VerifyText: Please enter your password.'
$ws.Range("L26").Value = 'VerifyText: Please enter your confirm password.'
$ws.Range("L27").Value = 'SetText: Akash'
$ws.Range("L28").Value = 'SetText: sangal'
$ws.Range("L29").Value = 'ClearText: null'
$ws.Range("L30").Value = 'SetText: 123456'
$ws.Range("L31").Value = 'SetText: 123456'
$ws.Range("L32").Value = 'Click: null'
$ws.Range("L33").Value = 'VerifyText: Please enter your email address.'
$ws.Range("L34").Value = 'SetText: Akash'
$ws.Range("L35").Value = 'SetText: sangal'
$ws.Range("L36").Value = 'SetText: komal@bravvura.in'
$ws.Range("L37").Value = 'ClearText: null'
$ws.Range("L38").Value = 'ClearText: null'
$ws.Range("L39").Value = 'Click: null'
$ws.Range("L40").Value = 'VerifyText: Please enter your password.'
$ws.Range("L41").Value = 'VerifyText: Please enter your confirm password.'
$ws.Range("L42").Value = 'SetText: Akash'
$ws.Range("L43").Value = 'SetText: sangal'
$ws.Range("L44").Value = 'SetText: komal@'
$ws.Range("L45").Value = 'SetText: 123456'
$ws.Range("L46").Value = 'SetText: 123456'
$ws.Range("L47").Value = 'Click: null'
$ws.Range("L48").Value = 'VerifyText: Please enter a valid email address (Ex: johndoe@domain.com).'
$ws.Range("L49").Value = 'SetText: Akash'
$ws.Range("L50").Value = 'SetText: sangal'
$ws.Range("L51").Value = 'SetText: komal@bravvura.in'
$ws.Range("L52").Value = 'SetText: 123456'
$ws.Range("L53").Value = 'SetText: 1234567'
$ws.Range("L54").Value = 'Click: null'
$ws.Range("L55").Value = 'VerifyText: Please make sure your passwords match.'
$ws.Range("L56").Value = 'SetText: Akash'
$ws.Range("L57").Value = 'SetText: sangal'
$ws.Range("L58").Value = 'SetText: komal@'
$ws.Range("L59").Value = 'SetText: 123456'
$ws.Range("L60").Value = 'SetText: 1234567'
$ws.Range("L61").Value = 'Click: null'
$ws.Range("L62").Value = 'VerifyText: Please make sure your passwords match.'
$ws.Range("L63").Value = 'VerifyText: Please enter a valid email address (Ex: johndoe@domain.com).'
$ws.Range("L64").Value = 'SetText: Akash'
$ws.Range("L65").Value = 'SetText: sangal'
$ws.Range("L66").Value = 'SetText: komal@d.com'
$ws.Range("L67").Value = 'SetText: 123'
$ws.Range("L68").Value = 'SetText: 123'
$ws.Range("L69").Value = 'Click: null'
$ws.Range("L70").Value = 'VerifyText: Please enter 6 or more characters. Leading and trailing spaces will be ignored.'
$ws.Range("L71").Value = 'Click: null'
$ws.Range("L72").Value = 'VerifyElement: null'
$ws.Range("L73").Value = 'Click: null'
$ws.Range("L74").Value = 'Click: null'
$ws.Range("L75").Value = 'VerifyElement: null'
$ws.Range("L76").Value = 'Click: null'
$ws.Range("L77").Value = 'SetText: Akash'
$ws.Range("L78").Value = 'SetText: sangal'
$ws.Range("L79").Value = 'SetText: Randomemailid'
$ws.Range("L80").Value = 'SetText: 123456'
$ws.Range("L81").Value = 'SetText: 123456'
$ws.Range("L82").Value = 'Click: null'
$ws.Range("L83").Value = 'Wait: 6000'
$ws.Range("L84").Value = 'VerifyText: Akash Sangal'
$ws.Range("L85").Value = 'Click: null'
$ws.Range("L86").Value = 'Click: null'
$ws.Range("L87").Value = 'Wait: 6000'
$ws.Range("L88").Value = 'VerifyElement: null'
$ws.Range("L89").Value = 'Click: null'
$ws.Range("L90").Value = 'SetText: Randomemailid'
$ws.Range("L91").Value = 'SetText: 123456'
$ws.Range("L92").Value = 'Click: null'
$ws.Range("L93").Value = 'VerifyText: Akash Sangal'
$ws.Range("L94").Value = 'Click: null'
$ws.Range("L95").Value = 'Click: null'
$ws.Range("L96").Value = 'Wait: 6000'
$ws.Range("L97").Value = 'VerifyElement: null'
$ws.Range("L98").Value = 'Click: null'
$ws.Range("L99").Value = 'VerifyElement: null'
$ws.Range("L100").Value = 'VerifyElement: null'
$ws.Range("L101").Value = 'Click: null'
$ws.Range("L102").Value = 'SetText: Akash'
$ws.Range("L103").Value = 'SetText: sangal'
$ws.Range("L104").Value = 'SetText: Randomemailid'
$ws.Range("L105").Value = 'SetText: 123456'
$ws.Range("L106").Value = 'SetText: 123456'
$ws.Range("L107").Value = 'Click: null'
$ws.Range("L108").Value = 'VerifyText: A customer with the same email already exists in an associated website.'
$ws.Range("L109").Value = 'Click: null'
$ws.Range("L110").Value = 'Click: null'
$ws.Range("L111").Value = 'Click: null'
$ws.Range("L112").Value = 'Click: null'
$ws.Range("L113").Value = 'VerifyNoElement: null'
